$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old row 14 ("Admin Add Product Category"). This shifts rows
#    15-25 up by one (row 25 disappears, matching the new 23-item list),
#    and also causes Excel to prune the now-unused shared string on save.
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. The old row 13 ("Admin Ban user by username") becomes the merged /
#    renamed entry "Admin Ban user when on profile".
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "Admin Ban user when on profile"

# ---------------------------------------------------------------------------
# 3. Fix up the sequential index numbers in column A for every row that
#    shifted up (rows 14-24 now hold what used to be rows 15-25, so their
#    index needs to be decremented by one).
# ---------------------------------------------------------------------------
for ($r = 14; $r -le 24; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------------
# 4. "Add product" (row 10) is now marked Solved -> switch its status cell
#    from the red fill style to the green fill style (copy format from any
#    existing "Solved" cell).
# ---------------------------------------------------------------------------
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# ---------------------------------------------------------------------------
# 5. "See categories and filter by them" (row 14) and "Search Bar with
#    search filter + categories" (row 17) get highlighted with a new yellow
#    status style (keeps the same border, default font, just a new fill).
# ---------------------------------------------------------------------------
$ws.Range("C11").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C14").Interior.Color = 65535           # RGB(FFFF00)

$ws.Range("C11").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C17").Interior.Color = 65535           # RGB(FFFF00)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 6. Update the sheet view: scroll back to the top-left corner and change
#    the active selection to C20.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C20").Select() | Out-Null
